# Apply the changes described by the commit "wrote naive GPU matrix chapter":
#  - Add a new log entry row (row 78) on Sheet1: date 2013-05-01 and the note
#    "Wrote naive GPU matrix chapter"
#  - Remove the unused, empty Sheet2 and Sheet3 worksheets

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row at the bottom of the activity log
$ws.Range("A78").Value = (Get-Date -Year 2013 -Month 5 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B78").Value = "Wrote naive GPU matrix chapter"

# Remove the now-unused empty sheets
[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()

$excel.DisplayAlerts = $true
